$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated fitted_prob / fitted_length / region values for existing rows (2-99).
# Each entry: row, column index, type ("N" = numeric, "S" = string), value
$updates = @(
    @(2, 1, "N", "5.76009798229831e-006"),
    @(3, 1, "N", "0.0002049382229491399"),
    @(4, 1, "N", "0.002090612404522586"),
    @(5, 1, "N", "0.00991646342110409"),
    @(6, 1, "N", "0.02874319209472026"),
    @(7, 1, "N", "0.05919016053999403"),
    @(8, 1, "N", "0.0943801149291129"),
    @(9, 1, "N", "0.1228489444620447"),
    @(10, 1, "N", "0.1367159256098504"),
    @(11, 1, "N", "0.1337776693127749"),
    @(12, 1, "N", "0.1176133185591873"),
    @(13, 1, "N", "0.09474360742810109"),
    @(14, 1, "N", "0.07069216672148865"),
    @(15, 1, "N", "0.049085736026836"),
    @(16, 1, "N", "0.03226170247274943"),
    @(17, 1, "N", "0.02039711328142171"),
    @(18, 1, "N", "0.0122587011585376"),
    @(19, 1, "N", "0.006953347753789159"),
    @(20, 1, "N", "0.003888369301103059"),
    @(21, 1, "N", "0.002060296099352595"),
    @(22, 1, "N", "0.001103816671239376"),
    @(23, 1, "N", "0.0005350827862503429"),
    @(24, 1, "N", "0.0002883080621666154"),
    @(25, 1, "N", "0.0001264189925588629"),
    @(26, 1, "N", "6.245158865018167e-005"),
    @(27, 1, "N", "3.061946822169101e-005"),
    @(28, 1, "N", "1.394550037819591e-005"),
    @(29, 1, "N", "6.972750189097953e-006"),
    @(30, 1, "N", "2.728467465299199e-006"),
    @(31, 1, "N", "9.09489155099733e-007"),
    @(32, 1, "N", "6.063261033998221e-007"),
    @(32, 2, "N", "32"),
    @(33, 1, "N", "8.953487705786961e-008"),
    @(33, 2, "N", "3"),
    @(33, 3, "S", "NJ"),
    @(34, 1, "N", "4.208139221719872e-006"),
    @(35, 1, "N", "9.633952771426771e-005"),
    @(36, 1, "N", "0.0008909615616028606"),
    @(37, 1, "N", "0.004885739171293829"),
    @(38, 1, "N", "0.01656780225541937"),
    @(39, 1, "N", "0.04061257255906437"),
    @(40, 1, "N", "0.07491687582013948"),
    @(41, 1, "N", "0.1114943800748368"),
    @(42, 1, "N", "0.1375617432512191"),
    @(43, 1, "N", "0.1458144414742741"),
    @(44, 1, "N", "0.1348603864755061"),
    @(45, 1, "N", "0.1114221254290511"),
    @(46, 1, "N", "0.08329966821955957"),
    @(47, 1, "N", "0.05695098645946147"),
    @(48, 1, "N", "0.0358907717476635"),
    @(49, 1, "N", "0.02122173563001037"),
    @(50, 1, "N", "0.0117123258725711"),
    @(51, 1, "N", "0.006141734426661624"),
    @(52, 1, "N", "0.003052602098411007"),
    @(53, 1, "N", "0.001457448728748002"),
    @(54, 1, "N", "0.0006530673932601009"),
    @(55, 1, "N", "0.0002867802112163564"),
    @(56, 1, "N", "0.0001258860371433647"),
    @(57, 1, "N", "5.094534504592781e-005"),
    @(58, 1, "N", "1.835464979686327e-005"),
    @(59, 1, "N", "6.356976271108742e-006"),
    @(60, 1, "N", "2.775581188793958e-006"),
    @(61, 1, "N", "6.267441394050873e-007"),
    @(62, 1, "N", "1.790697541157392e-007"),
    @(63, 1, "N", "8.953487705786961e-008"),
    @(63, 2, "N", "34"),
    @(64, 1, "N", "2.965874939533225e-007"),
    @(66, 1, "N", "0.0001474039844948013"),
    @(67, 1, "N", "0.001041318691270115"),
    @(68, 1, "N", "0.004535020507542269"),
    @(69, 1, "N", "0.01364302472185283"),
    @(70, 1, "N", "0.03098756023074112"),
    @(71, 1, "N", "0.05605711246963562"),
    @(72, 1, "N", "0.08436006156760925"),
    @(73, 1, "N", "0.1088147879315385"),
    @(74, 1, "N", "0.1233106994235031"),
    @(75, 1, "N", "0.1251313511863846"),
    @(76, 1, "N", "0.1151857838891498"),
    @(77, 1, "N", "0.09758984104788712"),
    @(78, 1, "N", "0.07706084561642201"),
    @(79, 1, "N", "0.05701469462511291"),
    @(80, 1, "N", "0.0396748056536299"),
    @(81, 1, "N", "0.02619332225348364"),
    @(82, 1, "N", "0.0164980748011455"),
    @(83, 1, "N", "0.009992724708773325"),
    @(84, 1, "N", "0.005828636293668678"),
    @(85, 1, "N", "0.003260485183526858"),
    @(86, 1, "N", "0.001792673675951866"),
    @(87, 1, "N", "0.0009433459557675343"),
    @(88, 1, "N", "0.000476121790293067"),
    @(89, 1, "N", "0.0002325245952594048"),
    @(90, 1, "N", "0.0001138895976780758"),
    @(91, 1, "N", "5.743911132896012e-005"),
    @(92, 1, "N", "2.392472451223468e-005"),
    @(93, 1, "N", "1.057828728433517e-005"),
    @(94, 1, "N", "5.338574891159804e-006"),
    @(95, 1, "N", "3.262462433486547e-006"),
    @(96, 1, "N", "1.680662465735494e-006"),
    @(97, 1, "N", "4.943124899222041e-007"),
    @(99, 1, "N", "1.977249959688816e-007")
)

foreach ($item in $updates) {
    $r = $item[0]
    $c = $item[1]
    $kind = $item[2]
    $raw = $item[3]
    if ($kind -eq "N") {
        $ws.Cells.Item($r, $c).Value = [double]$raw
    } else {
        $ws.Cells.Item($r, $c).Value = $raw
    }
}

# Append two new rows (100, 101) of NO region data for fitted_length 39 and 40.
$ws.Cells.Item(100, 1).Value = [double]"9.886249798444082e-008"
$ws.Cells.Item(100, 2).Value = 39
$ws.Cells.Item(100, 3).Value = "NO"
$ws.Cells.Item(100, 4).Value = "y2"

$ws.Cells.Item(101, 1).Value = [double]"9.886249798444082e-008"
$ws.Cells.Item(101, 2).Value = 40
$ws.Cells.Item(101, 3).Value = "NO"
$ws.Cells.Item(101, 4).Value = "y2"
